# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    and populate it with the per-fund holding breakdown for 2022-Q1.
# 2. Insert a new leading row into "总计" summarising the 2022-Q1 quarter
#    (date / holding count / holding value), shifting the previous rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create & position the "2022-Q1" worksheet
# ---------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4sheet)
$newSheet.Name = "2022-Q1"

# Copy the header + column-A formatting (bold, bordered, centered cell
# style) from the "2021-Q4" sheet so the new sheet matches the look of
# its siblings. "2021-Q4" only has 5 rows of its own, so stretch the
# last (data) row's formatting down across the extra rows 6-8 this
# sheet needs.
$q4sheet.Range("A1:H5").Copy()
$newSheet.Range("A1:H5").PasteSpecial(-4122)  # xlPasteFormats
$q4sheet.Range("A5:H5").Copy()
$newSheet.Range("A6:H8").PasteSpecial(-4122)  # xlPasteFormats

$fundRows = @(
    @{ idx=0; code="005662"; name="嘉实金融精选股票A";             scale="10.58"; pos="90.89"; pct="9.43"; value="0.9977"; rank=2 },
    @{ idx=1; code="012671"; name="嘉实核心蓝筹混合型证券投资基金A"; scale="11.11"; pos="93.27"; pct="6.78"; value="0.7533"; rank=3 },
    @{ idx=2; code="005663"; name="嘉实金融精选股票C";             scale="3.00";  pos="90.89"; pct="9.43"; value="0.2829"; rank=2 },
    @{ idx=3; code="513690"; name="博时恒生港股通高股息率ETF";      scale="4.60";  pos="99.64"; pct="3.16"; value="0.1454"; rank=3 },
    @{ idx=4; code="009126"; name="嘉实基础产业优选股票A";          scale="2.17";  pos="90.79"; pct="6.31"; value="0.1369"; rank=3 },
    @{ idx=5; code="012672"; name="嘉实核心蓝筹混合型证券投资基金C"; scale="0.47";  pos="93.27"; pct="6.78"; value="0.0319"; rank=3 },
    @{ idx=6; code="009127"; name="嘉实基础产业优选股票C";          scale="0.11";  pos="90.79"; pct="6.31"; value="0.0069"; rank=3 }
)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

foreach ($row in $fundRows) {
    $r = 2 + $row.idx

    $newSheet.Range("A$r").Value = $row.idx

    # Columns B-G carry values that look numeric (fund codes with leading
    # zeros, decimal scale/percentage figures) but must stay text, exactly
    # like the source data -- pre-format as Text so Excel doesn't silently
    # coerce them to numbers (and strip the leading zeros on the codes).
    $textRange = $newSheet.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"

    $newSheet.Range("B$r").Value = $row.code
    $newSheet.Range("C$r").Value = $row.name
    $newSheet.Range("D$r").Value = $row.scale
    $newSheet.Range("E$r").Value = $row.pos
    $newSheet.Range("F$r").Value = $row.pct
    $newSheet.Range("G$r").Value = $row.value

    $newSheet.Range("H$r").Value = $row.rank
}

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert(-4121)  # xlShiftDown

# The row-insert leaves stray formatting on the new row; clear it and
# re-copy the column-A number style from the row underneath so it again
# matches its siblings (bold/border/centred, same as A3:A7).
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 2.35

# Column A is a plain (non-formula) running index, 0, 1, 2, ... -- the
# insert only shifted rows down, it didn't renumber them. Bump the
# index on every row that used to be above the insertion point.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}
